$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 853
$ws.Range("I19").Value = 175
$ws.Range("J19").Value = 1124.2
$ws.Range("K19").Value = 175
$ws.Range("L19").Value = 1124.2
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -1474.2
$ws.Range("H40").Value = 1783.6666
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 1820.4
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 1820.4
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2170.4
$ws.Range("H86").Value = 2350.111
$ws.Range("J86").Value = 2099.5
$ws.Range("L86").Value = 2099.5
$ws.Range("N86").Value = -4345.5
$ws.Range("H89").Value = 2350.111
$ws.Range("J89").Value = 2099.5
$ws.Range("L89").Value = 10497.5
$ws.Range("N89").Value = -21729.5
$ws.Range("H127").Value = 519.3125
$ws.Range("J127").Value = 1100
$ws.Range("L127").Value = 3300
$ws.Range("N127").Value = -13220
$ws.Range("H137").Value = 1076.8636
$ws.Range("I137").Value = 882.2941
$ws.Range("J137").Value = 1738.4
$ws.Range("K137").Value = 2646.8823
$ws.Range("L137").Value = 5215.200000000001
$ws.Range("M137").Value = -96.88229999999976
$ws.Range("N137").Value = -10315.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2704.4614
$ws.Range("I2").Value = 1506.4
$ws.Range("K2").Value = 1506.4
$ws.Range("M2").Value = -1393.4
$ws.Range("H74").Value = 568.087
$ws.Range("I74").Value = 524.93335
$ws.Range("J74").Value = 649
$ws.Range("K74").Value = 524.93335
$ws.Range("L74").Value = 649
$ws.Range("M74").Value = 349.06665
$ws.Range("N74").Value = -2397
$ws.Range("H77").Value = 568.087
$ws.Range("I77").Value = 524.93335
$ws.Range("J77").Value = 649
$ws.Range("K77").Value = 2624.66675
$ws.Range("L77").Value = 3245
$ws.Range("M77").Value = 1743.33325
$ws.Range("N77").Value = -11981
$ws.Range("H97").Value = 2065.2666
$ws.Range("I97").Value = 1838.1
$ws.Range("K97").Value = 1838.1
$ws.Range("M97").Value = -1342.1
$ws.Range("H102").Value = 3717.7222
$ws.Range("I102").Value = 2419.9092
$ws.Range("J102").Value = 5757.143
$ws.Range("K102").Value = 2419.9092
$ws.Range("L102").Value = 5757.143
$ws.Range("M102").Value = -797.9092000000001
$ws.Range("N102").Value = -9001.143
$ws.Range("H116").Value = 2704.4614
$ws.Range("I116").Value = 1506.4
$ws.Range("K116").Value = 1506.4
$ws.Range("M116").Value = 787.5999999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2704.4614
$ws.Range("I3").Value = 1506.4
$ws.Range("K3").Value = 1506.4
$ws.Range("M3").Value = -1392.4
$ws.Range("H40").Value = 18675.428
$ws.Range("J40").Value = 18675.428
$ws.Range("L40").Value = 18675.428
$ws.Range("N40").Value = -19205.428
$ws.Range("H42").Value = 75491.30499999999
$ws.Range("J42").Value = 75491.30499999999
$ws.Range("L42").Value = 75491.30499999999
$ws.Range("N42").Value = -76147.30499999999
$ws.Range("H96").Value = 14031.5
$ws.Range("I96").Value = 9255.1
$ws.Range("J96").Value = 25972.5
$ws.Range("K96").Value = 9255.1
$ws.Range("L96").Value = 25972.5
$ws.Range("M96").Value = -6509.1
$ws.Range("N96").Value = -31464.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2307.875
$ws.Range("I31").Value = 1603.7106
$ws.Range("K31").Value = 1603.7106
$ws.Range("M31").Value = -1308.7106
$ws.Range("H34").Value = 2307.875
$ws.Range("I34").Value = 1603.7106
$ws.Range("K34").Value = 1603.7106
$ws.Range("M34").Value = -1401.7106
$ws.Range("H58").Value = 1275.9048
$ws.Range("I58").Value = 821.1111
$ws.Range("J58").Value = 4004.6667
$ws.Range("K58").Value = 821.1111
$ws.Range("L58").Value = 4004.6667
$ws.Range("M58").Value = -618.1111
$ws.Range("N58").Value = -4410.6667
$ws.Range("H136").Value = 1275.9048
$ws.Range("I136").Value = 821.1111
$ws.Range("J136").Value = 4004.6667
$ws.Range("K136").Value = 2463.3333
$ws.Range("L136").Value = 12014.0001
$ws.Range("M136").Value = 86.66670000000022
$ws.Range("N136").Value = -17114.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 977.83606
$ws.Range("I113").Value = 800.875
$ws.Range("J113").Value = 1092.6216
$ws.Range("K113").Value = 2402.625
$ws.Range("L113").Value = 3277.8648
$ws.Range("M113").Value = -232.625
$ws.Range("N113").Value = -7617.864799999999
$ws.Range("H132").Value = 1992.5
$ws.Range("J132").Value = 2892.5
$ws.Range("L132").Value = 26032.5
$ws.Range("N132").Value = -31092.5
$ws.Range("H137").Value = 2205.1765
$ws.Range("I137").Value = 721
$ws.Range("J137").Value = 3874.875
$ws.Range("K137").Value = 2163
$ws.Range("L137").Value = 11624.625
$ws.Range("M137").Value = 2937
$ws.Range("N137").Value = -21824.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 626875
$ws.Range("I21").Value = 10000000
$ws.Range("K21").Value = 10000000
$ws.Range("M21").Value = -9999827
$ws.Range("H24").Value = 910727.25
$ws.Range("I24").Value = 10000000
$ws.Range("J24").Value = 1800
$ws.Range("K24").Value = 10000000
$ws.Range("L24").Value = 1800
$ws.Range("M24").Value = -9999827
$ws.Range("N24").Value = -2146
$ws.Range("H30").Value = 626875
$ws.Range("I30").Value = 10000000
$ws.Range("K30").Value = 10000000
$ws.Range("M30").Value = -9999895
$ws.Range("H97").Value = 2550.3572
$ws.Range("I97").Value = 2112
$ws.Range("J97").Value = 3134.8333
$ws.Range("K97").Value = 2112
$ws.Range("L97").Value = 3134.8333
$ws.Range("M97").Value = -1616
$ws.Range("N97").Value = -4126.8333
$ws.Range("H102").Value = 1980
$ws.Range("I102").Value = 1977.7778
$ws.Range("K102").Value = 1977.7778
$ws.Range("M102").Value = -355.7778000000001
$ws.Range("H122").Value = 1747
$ws.Range("I122").Value = 1764.2727
$ws.Range("J122").Value = 1723.25
$ws.Range("K122").Value = 5292.8181
$ws.Range("L122").Value = 5169.75
$ws.Range("M122").Value = -2842.8181
$ws.Range("N122").Value = -10069.75
$ws.Range("H126").Value = 23810490
$ws.Range("I126").Value = 980
$ws.Range("J126").Value = 33334294
$ws.Range("K126").Value = 2940
$ws.Range("L126").Value = 100002882
$ws.Range("M126").Value = -470
$ws.Range("N126").Value = -100007822

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2264.6365
$ws.Range("I7").Value = 1732.3334
$ws.Range("J7").Value = 4660
$ws.Range("K7").Value = 1732.3334
$ws.Range("L7").Value = 4660
$ws.Range("M7").Value = -1620.3334
$ws.Range("N7").Value = -4884
$ws.Range("H46").Value = 8155.2856
$ws.Range("I46").Value = 846.5
$ws.Range("J46").Value = 17900.334
$ws.Range("K46").Value = 846.5
$ws.Range("L46").Value = 17900.334
$ws.Range("M46").Value = -658.5
$ws.Range("N46").Value = -18276.334
$ws.Range("H61").Value = 2733.5
$ws.Range("I61").Value = 2167
$ws.Range("J61").Value = 3300
$ws.Range("K61").Value = 2167
$ws.Range("L61").Value = 3300
$ws.Range("M61").Value = -1965
$ws.Range("N61").Value = -3704
$ws.Range("H113").Value = 2733.5
$ws.Range("I113").Value = 2167
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2167
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = 3
$ws.Range("N113").Value = -7640
$ws.Range("H126").Value = 2264.6365
$ws.Range("I126").Value = 1732.3334
$ws.Range("J126").Value = 4660
$ws.Range("K126").Value = 5197.0002
$ws.Range("L126").Value = 13980
$ws.Range("M126").Value = -2727.0002
$ws.Range("N126").Value = -18920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 571.381
$ws.Range("I107").Value = 501.33334
$ws.Range("J107").Value = 746.5
$ws.Range("K107").Value = 1504.00002
$ws.Range("L107").Value = 2239.5
$ws.Range("M107").Value = 415.9999800000001
$ws.Range("N107").Value = -6079.5
$ws.Range("H113").Value = 468.28262
$ws.Range("I113").Value = 376.45456
$ws.Range("J113").Value = 701.38464
$ws.Range("K113").Value = 1129.36368
$ws.Range("L113").Value = 2104.15392
$ws.Range("M113").Value = 1040.63632
$ws.Range("N113").Value = -6444.15392
$ws.Range("H132").Value = 16234644
$ws.Range("I132").Value = 20835738
$ws.Range("K132").Value = 62507214
$ws.Range("M132").Value = -62504684
